# Apply updated odds values to rows 3 and 4 of the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 3.5
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 4
$ws.Range("L3").Value = 2.75
$ws.Range("Q3").Value = 1.93
$ws.Range("R3").Value = 1.88
$ws.Range("X3").Value = 17
$ws.Range("Y3").Value = 12
$ws.Range("AH3").Value = 8
$ws.Range("AI3").Value = 10
$ws.Range("AK3").Value = 19
$ws.Range("AL3").Value = 17
$ws.Range("AO3").Value = 19
$ws.Range("AP3").Value = 26
$ws.Range("AX3").Value = 11
$ws.Range("AZ3").Value = 41

# Row 4 updates
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.6
$ws.Range("J4").Value = 2.38
$ws.Range("K4").Value = 2.1
$ws.Range("X4").Value = 7
$ws.Range("AA4").Value = 15
$ws.Range("AE4").Value = 21
$ws.Range("AN4").Value = 3.5
$ws.Range("AO4").Value = 9

$wb.Save()
